$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price/Volume columns (D, E) so that
# values like "1.010" or "0.000008400" are not auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.442.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.699.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.19"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5478"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2744"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06458"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.08"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07706"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.692.05"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.552"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008400"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.81"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.478.94"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.953"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.259"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.011"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.16"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1328"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.904"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.86"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06252"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.86%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.608"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.615"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.042"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6178"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.413"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.772"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01643"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.119.42"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.149"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8795"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.33"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.851.35"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.61"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.243"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05289"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.160"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4304"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.02%  "
